# Auto update: 2025-11-29 04:04:56
# Update final score (K) and MACRO_SCORE (N) values on Sheet1 rows 2-4

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("K2").Value = 70
$ws.Range("N2").Value = 85.8724807945396

$ws.Range("K3").Value = 68.8
$ws.Range("N3").Value = 85.8724807945396

$ws.Range("K4").Value = 64.2
$ws.Range("N4").Value = 85.8724807945396
